$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("ITI")
$ws2.Select()
$win = $excel.ActiveWindow
$ws2.Range("A52:A54").Select()
$win.ScrollRow = 28
Write-Host $win.ScrollRow
Write-Host "done"
